# The shipped "styled.xlsx" pandas doc-guide example was regenerated
# (new getting-started docs build) and one of the sampled random values
# in the DataFrame came back as NaN instead of a float. pandas' Styler
# renders that cell blank (no red/negative styling) while every other
# cell keeps its value/formatting. Reproduce that: clear E5's content
# and drop it back to the plain (non-negative) black-text look instead
# of the red "negative number" look it had before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("E5")
$cell.ClearContents()
$cell.Font.Color = 0
